$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Update the "CSP182823" guarantor/CSP-number placeholder values to
# new, distinct CSP numbers on each of the five loan sheets.
# ------------------------------------------------------------------

$wsLoanDetail = $wb.Worksheets.Item("LoanDetail")
$wsLoanDetail.Cells.Item(38, 2).Value = "CSP186606"
$wsLoanDetail.Cells.Item(39, 2).Value = "CSP186606"

$wsWeekendLoan = $wb.Worksheets.Item("WeekendLoan")
$wsWeekendLoan.Cells.Item(17, 2).Value = "CSP186625"
$wsWeekendLoan.Cells.Item(19, 2).Value = "CSP186625"

$wsShortTermLoan = $wb.Worksheets.Item("ShortTermLoan")
$wsShortTermLoan.Cells.Item(17, 2).Value = "CSP186834"
$wsShortTermLoan.Cells.Item(19, 2).Value = "CSP186834"

$wsOtherProductLoan = $wb.Worksheets.Item("OtherProductLoan")
$wsOtherProductLoan.Cells.Item(35, 2).Value = "CSP186912"

$wsCMSLoan = $wb.Worksheets.Item("CMSLoan")
$wsCMSLoan.Cells.Item(17, 2).Value = "CSP200048"
$wsCMSLoan.Cells.Item(19, 2).Value = "CSP200048"

# ------------------------------------------------------------------
# Update the saved cursor/selection position on a few sheets to match
# where the author last left off editing.
# ------------------------------------------------------------------

$wsWeekendLoan.Activate()
$wsWeekendLoan.Range("B17").Select()

$wsShortTermLoan.Activate()
$wsShortTermLoan.Range("B19").Select()

$wsOtherProductLoan.Activate()
$wsOtherProductLoan.Range("B35").Select()

# Restore CMSLoan as the active/selected tab (matches the saved
# workbook view), with its selection updated as well.
$wsCMSLoan.Activate()
$wsCMSLoan.Range("B17").Select()
